# Update stock data: quantity column C for rows 217-232 (250 -> 265)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 217; $r -le 232; $r++) {
    $ws.Cells.Item($r, 3).Value = 265
}

# Reflect the selection/view change recorded in the edit: the active
# cell/selection moved to C217:C232 (the block that was just updated).
$ws.Range("C217:C232").Select()
